$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Colors (Interior.Color values matching the existing theme fills) ---
$Gray   = 12566463   # header row fill
$Yellow = 65535      # "yellow" note fill
$Green  = 5296274    # "green" note fill

# NOTE: the order in which brand-new (never-before-seen) string values are
# assigned determines the order they're appended to the shared string table,
# so cell values are set in a specific sequence below (matching the order
# the author actually typed them in) rather than strict row order.

# Row 2 - header (A2 text unchanged, B2 note changes)
$ws.Cells.Item(2, 1).Value = "Things I don't like"
$ws.Cells.Item(2, 2).Value = "…but I'll worry about later"

# Row 3
$ws.Cells.Item(3, 1).Value = 'We need to figure out how "HasMoved" gets set'
$ws.Cells.Item(3, 2).Value = "Production code not in its own project"

# Row 4
$ws.Cells.Item(4, 1).Value = "I wish Rook.GetMovesFrom were more Linq-y"
$ws.Cells.Item(4, 2).Value = "Pawn doesn't have en passant move capability"

# Row 5
$ws.Cells.Item(5, 1).Value = "We've got Rook, so why not Bishop"
$ws.Cells.Item(5, 2).Value = "Pawn doesn't have diagonal capture capability"

# Rows 8-11 (typed before rows 6-7 were inserted, per shared-string order)
$ws.Cells.Item(8, 1).Value = "And queen?"
$ws.Cells.Item(9, 1).Value = "And king?"
$ws.Cells.Item(10, 1).Value = "And maybe knight"
$ws.Cells.Item(11, 1).Value = "I'm going to have the same default board size in a lot of tests"

# Rows 6-7 (B6/B7 are empty but still-styled cells)
$ws.Cells.Item(6, 1).Value = "At some point, additional tests on bishop.GetMovesFrom"
$ws.Cells.Item(6, 2).Value = ""
$ws.Cells.Item(7, 1).Value = "Refactor ugly Bishop.GetMovesFrom"
$ws.Cells.Item(7, 2).Value = ""

# --- Formatting (fills / bold) ---
# Row 2 already carries the bold/gray header style from the source
# workbook; it is left untouched so the existing style is reused as-is.

$ws.Cells.Item(3, 1).Interior.Color = $Yellow
$ws.Cells.Item(3, 2).Interior.Color = $Yellow

$ws.Cells.Item(4, 1).Interior.Color = $Green
$ws.Cells.Item(4, 2).Interior.Color = $Yellow

$ws.Cells.Item(5, 1).Interior.Color = $Green
$ws.Cells.Item(5, 2).Interior.Color = $Yellow

$ws.Cells.Item(6, 1).Interior.Color = $Yellow
$ws.Cells.Item(6, 2).Interior.Color = $Yellow

$ws.Cells.Item(7, 1).Interior.Color = $Yellow
$ws.Cells.Item(7, 2).Interior.Color = $Yellow

$ws.Cells.Item(8, 1).Interior.Color = $Yellow
$ws.Cells.Item(9, 1).Interior.Color = $Yellow
$ws.Cells.Item(10, 1).Interior.Color = $Yellow
$ws.Cells.Item(11, 1).Interior.Color = $Yellow

# --- Column A width: re-fit to the new (longer) content ---
$ws.Columns.Item(1).AutoFit()

# --- Selection moves to A8 ---
[void]$ws.Range("A8").Select()
